$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Add()
$ws3.Name = "RulesTest"
$ws3.Range("B17").Value = "1"
Write-Host "done"
